# "Generate Report for Handback" -- the handback for the two files has
# completed (zh-cn first, de-de a little later). Update status text,
# stamp handback datetimes, link the "Latest Target File" cells back to
# the same source doc pages, and widen the columns that now hold the
# longer status text / long filenames.

$wb = $excel.ActiveWorkbook

$newStatus   = "Handed back: in sync with en-US"
$zhHandback  = "2016-10-24 10:03:09"
$deHandback  = "2016-10-24 10:03:27"

$url4fef = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/45a387e26dbfe977afe7aa5fab247af75755b490/e2e/4fef9958-a67e-4b2d-8e05-464a4eac7091.md"
$url7ebb = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/45a387e26dbfe977afe7aa5fab247af75755b490/e2e/7ebb8ac3-68bf-4c58-9be0-69ba0332c46c.md"

$name4fef = "4fef9958-a67e-4b2d-8e05-464a4eac7091.md"
$name7ebb = "7ebb8ac3-68bf-4c58-9be0-69ba0332c46c.md"

$xlf4fefZh = "4fef9958-a67e-4b2d-8e05-464a4eac7091.f5a96d2708f49c7d0ae1272cbe2de4d71548c3a9.zh-cn.xlf"
$xlf7ebbZh = "7ebb8ac3-68bf-4c58-9be0-69ba0332c46c.635c2ee3c723ee3fd935d2cb37dc911d1eb61c73.zh-cn.xlf"
$xlf4fefDe = "4fef9958-a67e-4b2d-8e05-464a4eac7091.f5a96d2708f49c7d0ae1272cbe2de4d71548c3a9.de-de.xlf"
$xlf7ebbDe = "7ebb8ac3-68bf-4c58-9be0-69ba0332c46c.635c2ee3c723ee3fd935d2cb37dc911d1eb61c73.de-de.xlf"

# ---------------------------------------------------------------------
# Overview sheet: "Ready for handoff" -> "Handed back: in sync with en-US"
# for both locale columns, and widen those two columns to fit the text.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.9777050018311
$overview.Columns.Item(6).ColumnWidth = 29.9777050018311

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Status column picks up the same new text.
$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus
$zh.Columns.Item(3).ColumnWidth = 29.9777050018311

# Latest Target File (I) / Latest Handback File (J) are now populated;
# I links back to the source doc the same way column A does.
$zh.Range("I2").Value = $name4fef
$zh.Range("J2").Value = $xlf4fefZh
$zh.Range("I3").Value = $name7ebb
$zh.Range("J3").Value = $xlf7ebbZh

$zh.Hyperlinks.Add($zh.Range("I2"), $url4fef, "", "", $name4fef) | Out-Null
$zh.Hyperlinks.Add($zh.Range("I3"), $url7ebb, "", "", $name7ebb) | Out-Null

# Latest Handback DateTime (K) -- handback finished for zh-cn.
$zh.Range("K2").Value = $zhHandback
$zh.Range("K3").Value = $zhHandback

$zh.Columns.Item(9).ColumnWidth = 40
$zh.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------
# de-de sheet (mirrors zh-cn, but handback lands a little later with a
# different timestamp)
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus
$de.Columns.Item(3).ColumnWidth = 29.9777050018311

$de.Range("I2").Value = $name4fef
$de.Range("J2").Value = $xlf4fefDe
$de.Range("I3").Value = $name7ebb
$de.Range("J3").Value = $xlf7ebbDe

$de.Hyperlinks.Add($de.Range("I2"), $url4fef, "", "", $name4fef) | Out-Null
$de.Hyperlinks.Add($de.Range("I3"), $url7ebb, "", "", $name7ebb) | Out-Null

$de.Range("K2").Value = $deHandback
$de.Range("K3").Value = $deHandback

$de.Columns.Item(9).ColumnWidth = 40
$de.Columns.Item(10).ColumnWidth = 40
